$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.861.79'
$ws.Range('E2').Value = '  -1.23%  '
$ws.Range('D3').Value = '1.635.77'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  -0.37%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '215.17'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  -0.54%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '0.5021'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  -2.07%  '
$ws.Range('E7').Value = '  -0.36%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.2573'
$cell.Style = "Normal"
$ws.Range('E8').Value = '  -0.48%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.06418'
$cell.Style = "Normal"
$ws.Range('E9').Value = '  -0.18%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '19.65'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  -1.85%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.07697'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  -1.29%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.637.54'
$ws.Range('E12').Value = '  -1.24%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '4.244'
$cell.Style = "Normal"
$ws.Range('E13').Value = '  -1.30%  '
$ws.Range('D14').Value = '1.861.28'
$ws.Range('E14').Value = '  -1.30%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '0.5443'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  -1.95%  '
$ws.Range('D16').Value = '0.0₅7931'
$ws.Range('E16').Value = '  -1.40%  '
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').Value = '25.872.78'
$ws.Range('E18').Value = '  -1.30%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  -0.34%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '202.97'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  -3.56%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '4.329'
$cell.Style = "Normal"
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '9.945'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  -1.39%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '5.979'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  -0.68%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '1.003'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  -0.29%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '1.929'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  +11.25%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '140.71'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  -2.27%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '0.1145'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  -2.16%  '
$ws.Range('E28').Value = '  -0.60%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '6.709'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  -3.90%  '
$ws.Range('E30').Value = '  -0.72%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '0.05001'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  -2.37%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '3.261'
$cell.Style = "Normal"
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('E33').Value = '  -1.73%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '1.538'
$cell.Style = "Normal"
$ws.Range('E34').Value = '  -1.82%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '2.363'
$cell.Style = "Normal"
$ws.Range('E35').Value = '  -0.43%  '
$ws.Range('D36').Value = '1.170.58'
$ws.Range('E36').Value = '  +0.97%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '0.8943'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  -4.00%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '2.614'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  -5.13%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '0.5610'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  -1.83%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '0.01560'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  -2.13%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '2.557'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('E42').Value = '  -0.35%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '5.675'
$cell.Style = "Normal"
$ws.Range('E43').Value = '  +0.19%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '0.8071'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  -3.83%  '
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('D46').Value = '1.773.20'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('E47').Value = '  -0.53%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '0.4515'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('E49').Value = '  +0.05%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '54.74'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  -2.22%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.05083'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  +0.57%  '
